$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44559
$ws.Range("D3").Value = 44559
$ws.Range("D6").Value = 44574
$ws.Range("D7").Value = 44574
